# Auto-generated edit script: updates Price (D) and Volume(1h) (E) columns
# for the cryptos worksheet, matching the upstream GitHub Actions refresh.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 4).Value = "61.593.47"
$ws.Cells.Item(2, 5).Value = "  -4.24%  "
$ws.Cells.Item(3, 4).Value = "2.976.14"
$ws.Cells.Item(3, 5).Value = "  -5.13%  "
$ws.Cells.Item(4, 5).Value = "  +0.06%  "
$ws.Cells.Item(5, 4).Value = "'539.20"
$ws.Cells.Item(5, 5).Value = "  -5.65%  "
$ws.Cells.Item(6, 4).Value = "'150.79"
$ws.Cells.Item(6, 5).Value = "  -8.09%  "
$ws.Cells.Item(7, 5).Value = "  +0.12%  "
$ws.Cells.Item(8, 4).Value = "'0.567"
$ws.Cells.Item(8, 5).Value = "  -1.42%  "
$ws.Cells.Item(9, 4).Value = "2.986.32"
$ws.Cells.Item(9, 5).Value = "  -5.19%  "
$ws.Cells.Item(10, 5).Value = "  -3.77%  "
$ws.Cells.Item(11, 4).Value = "'6.14"
$ws.Cells.Item(11, 5).Value = "  -6.88%  "
$ws.Cells.Item(12, 5).Value = "  -4.10%  "
$ws.Cells.Item(13, 4).Value = "3.499.21"
$ws.Cells.Item(13, 5).Value = "  -5.18%  "
$ws.Cells.Item(14, 5).Value = "  -1.50%  "
$ws.Cells.Item(15, 4).Value = "61.648.95"
$ws.Cells.Item(15, 5).Value = "  -4.07%  "
$ws.Cells.Item(16, 4).Value = "'23.61"
$ws.Cells.Item(16, 5).Value = "  -5.74%  "
$ws.Cells.Item(17, 4).Value = "2.984.20"
$ws.Cells.Item(17, 5).Value = "  -4.84%  "
$ws.Cells.Item(18, 5).Value = "  -5.76%  "
$ws.Cells.Item(19, 5).Value = "  -1.93%  "
$ws.Cells.Item(20, 5).Value = "  -3.83%  "
$ws.Cells.Item(21, 4).Value = "'379.98"
$ws.Cells.Item(21, 5).Value = "  -6.39%  "
$ws.Cells.Item(22, 4).Value = "'6.68"
$ws.Cells.Item(22, 5).Value = "  -5.67%  "
$ws.Cells.Item(23, 5).Value = "  +0.07%  "
$ws.Cells.Item(24, 5).Value = "  -3.56%  "
$ws.Cells.Item(25, 4).Value = "'65.79"
$ws.Cells.Item(25, 5).Value = "  -4.58%  "
$ws.Cells.Item(26, 4).Value = "'0.471"
$ws.Cells.Item(26, 5).Value = "  -2.90%  "
$ws.Cells.Item(27, 4).Value = "3.101.89"
$ws.Cells.Item(27, 5).Value = "  -5.29%  "
$ws.Cells.Item(28, 4).Value = "'0.189"
$ws.Cells.Item(28, 5).Value = "  -3.19%  "
$ws.Cells.Item(29, 4).Value = "'0.997"
$ws.Cells.Item(29, 5).Value = "  +0.01%  "
$ws.Cells.Item(30, 4).Value = "0.0₃0933"
$ws.Cells.Item(30, 5).Value = "  -8.41%  "
$ws.Cells.Item(31, 4).Value = "'8.16"
$ws.Cells.Item(31, 5).Value = "  -8.11%  "
$ws.Cells.Item(32, 5).Value = "  +0.01%  "
$ws.Cells.Item(33, 5).Value = "  -4.75%  "
$ws.Cells.Item(34, 4).Value = "'20.45"
$ws.Cells.Item(34, 5).Value = "  -3.62%  "
$ws.Cells.Item(35, 4).Value = "'158.80"
$ws.Cells.Item(35, 5).Value = "  -1.88%  "
$ws.Cells.Item(36, 4).Value = "'4.56"
$ws.Cells.Item(36, 5).Value = "  -6.38%  "
$ws.Cells.Item(37, 5).Value = "  -6.23%  "
$ws.Cells.Item(38, 4).Value = "'1.06"
$ws.Cells.Item(38, 5).Value = "  -4.66%  "
$ws.Cells.Item(39, 5).Value = "  -6.10%  "
$ws.Cells.Item(40, 5).Value = "  -8.47%  "
$ws.Cells.Item(41, 5).Value = "  -2.13%  "
$ws.Cells.Item(42, 2).Value = "Maker"
$ws.Cells.Item(42, 3).Value = "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
$ws.Cells.Item(42, 4).Value = "2.411.42"
$ws.Cells.Item(42, 5).Value = "  -8.82%  "
$ws.Cells.Item(43, 2).Value = "Filecoin"
$ws.Cells.Item(43, 3).Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Cells.Item(43, 4).Value = "'3.89"
$ws.Cells.Item(43, 5).Value = "  -4.38%  "
$ws.Cells.Item(44, 4).Value = "'22.02"
$ws.Cells.Item(44, 5).Value = "  -6.78%  "
$ws.Cells.Item(45, 5).Value = "  -2.71%  "
$ws.Cells.Item(46, 4).Value = "'0.0589"
$ws.Cells.Item(46, 5).Value = "  -3.88%  "
$ws.Cells.Item(47, 4).Value = "'5.09"
$ws.Cells.Item(47, 5).Value = "  -5.54%  "
$ws.Cells.Item(48, 5).Value = "  +0.08%  "
$ws.Cells.Item(49, 4).Value = "'0.0244"
$ws.Cells.Item(49, 5).Value = "  -3.69%  "
$ws.Cells.Item(50, 4).Value = "'0.0951"
$ws.Cells.Item(50, 5).Value = "  -2.66%  "
$ws.Cells.Item(51, 4).Value = "'19.72"
$ws.Cells.Item(51, 5).Value = "  -6.86%  "
